$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update item name (الاسم) in C7 (merged C7:G7)
$ws.Range("C7").Value = "CONA-ADIONE 10MG 30 CHEWABLE TAB."

# Update price (السعر) in N7 (merged N7:O7) -- force text storage to keep "54.00"
$cellN7 = $ws.Range("N7")
$fmtN7 = $cellN7.NumberFormat
$cellN7.NumberFormat = "@"
$cellN7.Value = "54.00"
$cellN7.NumberFormat = $fmtN7

# Update selling price (سعر البيع) in P7 -- force text storage to keep "54.0000"
$cellP7 = $ws.Range("P7")
$fmtP7 = $cellP7.NumberFormat
$cellP7.NumberFormat = "@"
$cellP7.Value = "54.0000"
$cellP7.NumberFormat = $fmtP7
